$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '68.672.73'
Set-TextValue 'E2' '  +0.28%  '
Set-TextValue 'D3' '2.443.08'
Set-TextValue 'E3' '  -0.50%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '560.17'
Set-TextValue 'E5' '  +0.94%  '
Set-TextValue 'D6' '162.96'
Set-TextValue 'E6' '  +1.20%  '
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'D8' '0.511'
Set-TextValue 'E8' '  +2.48%  '
Set-TextValue 'D9' '0.160'
Set-TextValue 'E9' '  +8.18%  '
Set-TextValue 'E10' '  +0.09%  '
Set-TextValue 'E11' '  -1.90%  '
Set-TextValue 'E12' '  +1.03%  '
Set-TextValue 'D13' '68.594.22'
Set-TextValue 'E13' '  +0.31%  '
Set-TextValue 'E14' '  +2.66%  '
Set-TextValue 'D15' '23.19'
Set-TextValue 'E15' '  +0.17%  '
Set-TextValue 'D16' '10.39'
Set-TextValue 'E16' '  -3.39%  '
Set-TextValue 'D17' '337.85'
Set-TextValue 'E17' '  -0.73%  '
Set-TextValue 'D18' '6.88'
Set-TextValue 'E18' '  -1.61%  '
Set-TextValue 'E19' '  +1.23%  '
Set-TextValue 'D20' '1.90'
Set-TextValue 'E20' '  +1.83%  '
Set-TextValue 'E21' '  -0.15%  '
Set-TextValue 'D22' '66.90'
Set-TextValue 'E22' '  +1.01%  '
Set-TextValue 'E23' '  +0.26%  '
Set-TextValue 'D24' '8.12'
Set-TextValue 'E24' '  +0.92%  '
Set-TextValue 'E25' '  +0.34%  '
Set-TextValue 'D26' '7.19'
Set-TextValue 'E26' '  +1.40%  '
Set-TextValue 'D27' '1.00'
Set-TextValue 'E27' '  +0.10%  '
Set-TextValue 'D28' '425.75'
Set-TextValue 'E28' '  -0.53%  '
Set-TextValue 'E29' '  +0.91%  '
Set-TextValue 'E30' '  +0.08%  '
Set-TextValue 'D31' '161.00'
Set-TextValue 'E31' '  +2.43%  '
Set-TextValue 'D32' '19.00'
Set-TextValue 'E32' '  -0.01%  '
Set-TextValue 'E33' '  -0.16%  '
Set-TextValue 'D34' '17.79'
Set-TextValue 'E34' '  +0.38%  '
Set-TextValue 'D35' '0.104'
Set-TextValue 'E35' '  -4.33%  '
Set-TextValue 'D36' '0.296'
Set-TextValue 'E36' '  -2.14%  '
Set-TextValue 'D37' '4.34'
Set-TextValue 'E37' '  -0.51%  '
Set-TextValue 'E38' '  +1.14%  '
Set-TextValue 'D39' '1.05'
Set-TextValue 'E39' '  -3.55%  '
Set-TextValue 'E40' '  -0.54%  '
Set-TextValue 'D41' '3.35'
Set-TextValue 'E41' '  +1.29%  '
Set-TextValue 'D42' '129.73'
Set-TextValue 'E42' '  -2.11%  '
Set-TextValue 'E43' '  +0.67%  '
Set-TextValue 'D44' '0.479'
Set-TextValue 'E44' '  +0.59%  '
Set-TextValue 'E45' '  +0.57%  '
Set-TextValue 'D46' '0.0918'
Set-TextValue 'E46' '  +1.37%  '
Set-TextValue 'E47' '  +1.05%  '
Set-TextValue 'E48' '  -3.26%  '
Set-TextValue 'D49' '4.88'
Set-TextValue 'E49' '  -4.16%  '
Set-TextValue 'D50' '16.62'
Set-TextValue 'E50' '  -0.94%  '
Set-TextValue 'E51' '  +3.24%  '
